# Preenche a planilha "automacoestemp" com os dados das automacoes
# (mensagem/confirmacao exibida ao usuario ao acionar cada automacao).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 1 - Ar condicionado
$ws.Range("A1").Value = "ar"
$ws.Range("B1").Value = "A/C"
$ws.Range("C1").Value = 16
$ws.Range("D1").Formula = "FALSE"
$ws.Range("F1").Formula = "FALSE"

# Linha 2 - Televisor
$ws.Range("A2").Value = "tv"
$ws.Range("B2").Value = "Televisor"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Formula = "FALSE"
$ws.Range("F2").Formula = "TRUE"

# Linha 3 - Lampada
$ws.Range("A3").Value = "lamp1"
$ws.Range("B3").Value = "Lâmpada"
$ws.Range("C3").Value = 0
$ws.Range("D3").Formula = "FALSE"
$ws.Range("F3").Formula = "TRUE"
